# Adds a "LoginUserBean" data-provider sheet (before the existing
# "AdminUserBean" sheet) and adds a new "searchUserDataProvider" row
# to the "AdminUserBean" sheet, matching the data-driven framework
# change described in the commit message.

$wb = $excel.ActiveWorkbook

# --- Create the new sheet; Worksheets.Add() inserts it right before
# --- the currently active sheet, i.e. as the new first tab.
$loginSheet = $wb.Worksheets.Add()
$loginSheet.Name = "LoginUserBean"

# Always look sheets up by name afterwards - an index/reference taken
# before the Add() call can point at the wrong sheet once the
# worksheet collection has shifted.
$adminSheet = $wb.Worksheets.Item("AdminUserBean")

# --- Populate LoginUserBean (header row + one data row) ---
$loginSheet.Range("A1").Value = "notation"
$loginSheet.Range("B1").Value = "username"
$loginSheet.Range("C1").Value = "password"
$loginSheet.Range("C2").Value = "admin123"
$loginSheet.Range("B2").Value = "Admin"
$loginSheet.Range("A2").Value = "validUserLoginDataProvider"

$loginSheet.Columns.Item(1).ColumnWidth = 23.166666666666668

$null = $loginSheet.Range("B8").Select()

# --- Add a new data row to AdminUserBean ---
$adminSheet.Range("B5").Value = "Admin"
$adminSheet.Range("E5").Value = "Enabled"
$adminSheet.Range("A5").Value = "searchUserDataProvider"
$adminSheet.Range("D5").Value = "Linda.Anderson"
$adminSheet.Range("C5").Value = "Linda Jane Anderson"

# Match A5's look to the other data-provider label cells (A2:A4) by
# copying their formatting instead of setting font properties by hand
# (which would create extra, unused font/style entries).
$null = $adminSheet.Range("A2").Copy()
$null = $adminSheet.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Resize columns C and D to fit the new, longer values.
$adminSheet.Columns.Item(3).ColumnWidth = 21
$adminSheet.Columns.Item(4).ColumnWidth = 12.833333333333334

# Remove the hyperlink that used to sit on D4.
$null = $adminSheet.Range("D4").Hyperlinks.Item(1).Delete()

$null = $adminSheet.Range("A8").Select()

# AdminUserBean ends up as the active tab.
$null = $adminSheet.Activate()
